$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arcline_payment")

# Row 15 ("SI-CH-133") was a single 256720.8 NOT_PAID invoice; it is now split
# into a paid portion (row 15, 100000, PAID) and the remaining balance
# (new row 16, 156720.8, NOT_PAID).
$ws.Cells.Item(15, 3).Value = 100000
$ws.Cells.Item(15, 4).Value = "PAID"

# Insert the new row for the remaining balance; this shifts the old rows
# 16-20 down to 17-21 and copies row 15's formatting into the new row 16.
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).Value = 45250
$ws.Cells.Item(16, 2).Value = "SI-CH-133"
$ws.Cells.Item(16, 3).Value = 156720.79999999999
$ws.Cells.Item(16, 4).Value = "NOT_PAID"

# Append two brand-new invoices at the bottom of the table.
$ws.Cells.Item(22, 1).Value = 45313
$ws.Cells.Item(21, 1).Copy()
$ws.Cells.Item(22, 1).PasteSpecial(-4122)
$ws.Cells.Item(22, 2).Value = "SI-KA-94"
$ws.Cells.Item(22, 3).Value = 147382
$ws.Cells.Item(22, 3).NumberFormat = "#,##0.00"
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4122)
$ws.Cells.Item(22, 4).Value = "NOT_PAID"

$ws.Cells.Item(23, 1).Value = 45314
$ws.Cells.Item(21, 1).Copy()
$ws.Cells.Item(23, 1).PasteSpecial(-4122)
$ws.Cells.Item(23, 2).Value = "SI-KA-98"
$ws.Cells.Item(23, 3).Value = 9322
$ws.Cells.Item(23, 3).NumberFormat = "#,##0.00"
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4122)
$ws.Cells.Item(23, 4).Value = "NOT_PAID"

# Move the selection/active cell and make this sheet the active tab
# (mirrors the workbook being saved while viewing Arcline_payment, D24).
$ws.Range("D24").Select()
